$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# ---------------------------------------------------------------------
# Helpers that work purely off paragraph *index* (1-based) so we never
# have to depend on matching exact (unicode-sensitive) source text.
# Paragraph boundaries are re-derived from the live text every time a
# helper runs, since earlier edits shift later offsets.
# ---------------------------------------------------------------------

function Get-ParaInfo($range, [int]$paraIndex) {
    $full = $range.Text
    $parts = $full.Split([char]13)
    $pos = 1
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($i -eq ($paraIndex - 1)) {
            return @{ Start = $pos; Length = $parts[$i].Length }
        }
        $pos += $parts[$i].Length + 1
    }
    return $null
}

function Set-ParaText($range, [int]$paraIndex, [string]$newText) {
    $info = Get-ParaInfo $range $paraIndex
    $sub = $range.Characters($info.Start, $info.Length)
    $sub.Text = $newText
}

function Remove-Para($range, [int]$paraIndex) {
    $info = Get-ParaInfo $range $paraIndex
    # include the trailing paragraph mark (+1) so the paragraph itself is
    # removed (merged away) rather than just left empty
    $len = $info.Length + 1
    $sub = $range.Characters($info.Start, $len)
    $sub.Delete()
}

$ndash = [char]0x2013

# 1) "Background - Started, but Relies on Level Design" -> "Background - "
$t1 = "Background " + $ndash + " "
Set-ParaText $tr 1 $t1

# 2) "HUD - ..." -> "HUD - Still discussing"
$t2 = "HUD " + $ndash + " Still discussing"
Set-ParaText $tr 2 $t2

# 3) "Troubleshoot our Ghost Enemy - ... - Not Done" -> "Troubleshoot our Ghost Enemy- Not Done"
$t3 = "Troubleshoot our Ghost Enemy" + $ndash + " Not Done"
Set-ParaText $tr 3 $t3

# 4) Repurpose "Sort out our Level Design - Still Discussing" (paragraph 4)
#    into the new "Begin working on Weapons, Bullets, etc. ..." paragraph.
$t4 = "Begin working on Weapons, Bullets, etc. " + $ndash + " Sprites Done! Stats still to come"
Set-ParaText $tr 4 $t4

# 5) Remove the three now-unwanted level-1 bullet paragraphs entirely:
#    "Figure out level dimensions", "Parkour Additions",
#    "Make things a level instead of an enemy test!" (all currently at
#    paragraph index 5, since each removal shifts the next one up).
Remove-Para $tr 5
Remove-Para $tr 5
Remove-Para $tr 5

# 6) Repurpose the old "Begin working on Weapons, Bullets, etc. - Discussing"
#    paragraph (now at index 5, right before "SFX - Done") into "Tilesets".
Set-ParaText $tr 5 "Tilesets"

# 7) bodyPr autofit: <a:normAutofit lnSpcReduction="10000"/> -> <a:normAutofit/>
$shp.TextFrame.AutoSize = 2
